$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "텐서 ≈ 디스코드 서버를 오픈했습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2023/10/19/%ed%85%90%ec%84%9c-%e2%89%88-%eb%94%94%ec%8a%a4%ec%bd%94%eb%93%9c-%ec%84%9c%eb%b2%84%eb%a5%bc-%ec%98%a4%ed%94%88%ed%96%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

$ws.Range("D27").Value = "루다 서버에서 루다의 개인화 메시지를 처리하는 방법"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/personal-message-with-annotation/"
